$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for case with 380 kV (Case_3_107)

# Row 2
$ws.Range("B2").Value = 11.52589072752721
$ws.Range("C2").Value = 9.327864458967298
$ws.Range("D2").Value = 5.635391328243341
$ws.Range("F2").Value = 26.95649362286134
$ws.Range("G2").Value = 3.653003339112604
$ws.Range("I2").Value = 24.68767083849919
$ws.Range("K2").Value = 8.616769373253556
$ws.Range("L2").Value = 11.0800575026089
$ws.Range("M2").Value = 14.39478059419321
$ws.Range("N2").Value = 19.89759458299134
$ws.Range("O2").Value = 24.21293480060377

# Row 3
$ws.Range("B3").Value = 11.29816485131697
$ws.Range("C3").Value = 9.305270544962809
$ws.Range("D3").Value = 5.579866193628829
$ws.Range("F3").Value = 26.9885128249812
$ws.Range("G3").Value = 3.654639568776583
$ws.Range("I3").Value = 24.76669235662604
$ws.Range("K3").Value = 8.445337834980247
$ws.Range("L3").Value = 11.0876034429498
$ws.Range("M3").Value = 14.36206386437603
$ws.Range("N3").Value = 19.9552531405293
$ws.Range("O3").Value = 24.27748628428855

# Row 4
$ws.Range("B4").Value = 11.15807040718629
$ws.Range("C4").Value = 9.291267010532584
$ws.Range("D4").Value = 5.544901278840164
$ws.Range("F4").Value = 27.01434453901285
$ws.Range("G4").Value = 3.65569846970947
$ws.Range("I4").Value = 24.81926132713584
$ws.Range("K4").Value = 8.339657132723897
$ws.Range("L4").Value = 11.09391644663053
$ws.Range("M4").Value = 14.34411936345656
$ws.Range("N4").Value = 19.99232466957786
$ws.Range("O4").Value = 24.32158993867182

# Row 5
$ws.Range("B5").Value = 11.10099301625842
$ws.Range("C5").Value = 9.285528696986566
$ws.Range("D5").Value = 5.530439247176644
$ws.Range("F5").Value = 27.02642187991433
$ws.Range("G5").Value = 3.656143663142726
$ws.Range("I5").Value = 24.84170157675592
$ws.Range("K5").Value = 8.296544477555001
$ws.Range("L5").Value = 11.09691203313343
$ws.Range("M5").Value = 14.33735138603994
$ws.Range("N5").Value = 20.00785249128133
$ws.Range("O5").Value = 24.34068507580309

# Row 6
$ws.Range("B6").Value = 11.09151847168222
$ws.Range("C6").Value = 9.284573968898574
$ws.Range("D6").Value = 5.528025111454759
$ws.Range("F6").Value = 27.02852092453921
$ws.Range("G6").Value = 3.656218414790145
$ws.Range("I6").Value = 24.84548923704871
$ws.Range("K6").Value = 8.289384571732931
$ws.Range("L6").Value = 11.09743501064643
$ws.Range("M6").Value = 14.33626061087448
$ws.Range("N6").Value = 20.01045632934239
$ws.Range("O6").Value = 24.34392357137574

# Row 7
$ws.Range("B7").Value = 11.15730048491792
$ws.Range("C7").Value = 9.29118974874075
$ws.Range("D7").Value = 5.5447070957426
$ws.Range("F7").Value = 27.01450114157081
$ws.Range("G7").Value = 3.655704418288181
$ws.Range("I7").Value = 24.81955984287831
$ws.Range("K7").Value = 8.339075811358814
$ws.Range("L7").Value = 11.0939551327905
$ws.Range("M7").Value = 14.34402587631054
$ws.Range("N7").Value = 19.99253237751809
$ws.Range("O7").Value = 24.32184291886957

# Row 8
$ws.Range("B8").Value = 11.44747431832478
$ws.Range("C8").Value = 9.32010189555886
$ws.Range("D8").Value = 5.616430736034845
$ws.Range("F8").Value = 26.96625230881523
$ws.Range("G8").Value = 3.653556276540449
$ws.Range("I8").Value = 24.71407682330308
$ws.Range("K8").Value = 8.557782067429146
$ws.Range("L8").Value = 11.08231109396254
$ws.Range("M8").Value = 14.38305835853915
$ws.Range("N8").Value = 19.91712959274892
$ws.Range("O8").Value = 24.23426397896544

# Row 9
$ws.Range("B9").Value = 12.01094996203841
$ws.Range("C9").Value = 9.375712585424688
$ws.Range("D9").Value = 5.749912695028844
$ws.Range("F9").Value = 26.92064009540286
$ws.Range("G9").Value = 3.649772338058537
$ws.Range("I9").Value = 24.53936835147679
$ws.Range("K9").Value = 8.980818118703283
$ws.Range("L9").Value = 11.07277594381147
$ws.Range("M9").Value = 14.47635225408601
$ws.Range("N9").Value = 19.78245136393314
$ws.Range("O9").Value = 24.09802936059841

# Row 10
$ws.Range("B10").Value = 12.41716661653795
$ws.Range("C10").Value = 9.415844783084653
$ws.Range("D10").Value = 5.843292735488368
$ws.Range("F10").Value = 26.91701493688197
$ws.Range("G10").Value = 3.647250905934647
$ws.Range("I10").Value = 24.43062262284961
$ws.Range("K10").Value = 9.284850873073635
$ws.Range("L10").Value = 11.07383337842891
$ws.Range("M10").Value = 14.55473829292453
$ws.Range("N10").Value = 19.69146242181452
$ws.Range("O10").Value = 24.01964767087606

# Row 11
$ws.Range("B11").Value = 12.59941509284374
$ws.Range("C11").Value = 9.433931038472352
$ws.Range("D11").Value = 5.884690565140672
$ws.Range("F11").Value = 26.92184485526328
$ws.Range("G11").Value = 3.646159439857997
$ws.Range("I11").Value = 24.38541363649436
$ws.Range("K11").Value = 9.421063317050846
$ws.Range("L11").Value = 11.0760534795724
$ws.Range("M11").Value = 14.59244892346871
$ws.Range("N11").Value = 19.65178045748404
$ws.Range("O11").Value = 23.98871546548772

# Row 12
$ws.Range("B12").Value = 12.66799550520874
$ws.Range("C12").Value = 9.440754159858486
$ws.Range("D12").Value = 5.900206045138713
$ws.Range("F12").Value = 26.92460342124934
$ws.Range("G12").Value = 3.645754075184873
$ws.Range("I12").Value = 24.36890709321074
$ws.Range("K12").Value = 9.472293699336573
$ws.Range("L12").Value = 11.07714303237736
$ws.Range("M12").Value = 14.60701616660779
$ws.Range("N12").Value = 19.63699848301474
$ws.Range("O12").Value = 23.97768231150702

# Row 13
$ws.Range("B13").Value = 12.65324578404079
$ws.Range("C13").Value = 9.439285844554258
$ws.Range("D13").Value = 5.89687175433505
$ws.Range("F13").Value = 26.92396801206593
$ws.Range("G13").Value = 3.645841024744451
$ws.Range("I13").Value = 24.37243479565937
$ws.Range("K13").Value = 9.461276652304736
$ws.Range("L13").Value = 11.07689733093213
$ws.Range("M13").Value = 14.60386621618657
$ws.Range("N13").Value = 19.6401711786102
$ws.Range("O13").Value = 23.98002823368709

# Row 14
$ws.Range("B14").Value = 12.60506629834265
$ws.Range("C14").Value = 9.434492904302253
$ws.Range("D14").Value = 5.885970290206194
$ws.Range("F14").Value = 26.92205318842012
$ws.Range("G14").Value = 3.646125931150015
$ws.Range("I14").Value = 24.38404334241691
$ws.Range("K14").Value = 9.425285363640528
$ws.Range("L14").Value = 11.07613813790749
$ws.Range("M14").Value = 14.59364167406574
$ws.Range("N14").Value = 19.65055943705923
$ws.Range("O14").Value = 23.98779412438246

# Row 15
$ws.Range("B15").Value = 12.57549662696183
$ws.Range("C15").Value = 9.431553697933763
$ws.Range("D15").Value = 5.879271703931242
$ws.Range("F15").Value = 26.92100128853612
$ws.Range("G15").Value = 3.646301478743407
$ws.Range("I15").Value = 24.39123377429829
$ws.Range("K15").Value = 9.403192628539635
$ws.Range("L15").Value = 11.07570547917133
$ws.Range("M15").Value = 14.58741598262961
$ws.Range("N15").Value = 19.65695438375715
$ws.Range("O15").Value = 23.99263956316922

# Row 16
$ws.Range("B16").Value = 12.40519923194401
$ws.Range("C16").Value = 9.414659270904055
$ws.Range("D16").Value = 5.840565058974771
$ws.Range("F16").Value = 26.91682950153741
$ws.Range("G16").Value = 3.647323350315236
$ws.Range("I16").Value = 24.4336629046404
$ws.Range("K16").Value = 9.275902661698126
$ws.Range("L16").Value = 11.07372318201724
$ws.Range("M16").Value = 14.55231443092822
$ws.Range("N16").Value = 19.69409003568099
$ws.Range("O16").Value = 24.02176430185551

# Row 17
$ws.Range("B17").Value = 12.30002693388629
$ws.Range("C17").Value = 9.404250607493369
$ws.Range("D17").Value = 5.816538892960301
$ws.Range("F17").Value = 26.91592840674458
$ws.Range("G17").Value = 3.647964435235174
$ws.Range("I17").Value = 24.46078327517539
$ws.Range("K17").Value = 9.19724217283461
$ws.Range("L17").Value = 11.07295169014272
$ws.Range("M17").Value = 14.53130079636068
$ws.Range("N17").Value = 19.71730859216539
$ws.Range("O17").Value = 24.04084199641327

# Row 18
$ws.Range("B18").Value = 12.23929918588959
$ws.Range("C18").Value = 9.398247888569102
$ws.Range("D18").Value = 5.802618227023895
$ws.Range("F18").Value = 26.91602012592732
$ws.Range("G18").Value = 3.648338400809354
$ws.Range("I18").Value = 24.47678313437756
$ws.Range("K18").Value = 9.151804440538081
$ws.Range("L18").Value = 11.07267171906311
$ws.Range("M18").Value = 14.51940800064058
$ws.Range("N18").Value = 19.73082424885085
$ws.Range("O18").Value = 24.05225958336226

# Row 19
$ws.Range("B19").Value = 12.21869950730589
$ws.Range("C19").Value = 9.396212774172419
$ws.Range("D19").Value = 5.797887674825325
$ws.Range("F19").Value = 26.91615600195897
$ws.Range("G19").Value = 3.648465918619253
$ws.Range("I19").Value = 24.48226925688532
$ws.Range("K19").Value = 9.136388184169464
$ws.Range("L19").Value = 11.07260508880244
$ws.Range("M19").Value = 14.51541481180697
$ws.Range("N19").Value = 19.73542809559708
$ws.Range("O19").Value = 24.05620171356929

# Row 20
$ws.Range("B20").Value = 12.31124757301226
$ws.Range("C20").Value = 9.405360285628523
$ws.Range("D20").Value = 5.819107059491035
$ws.Range("F20").Value = 26.9159612095639
$ws.Range("G20").Value = 3.647895649589057
$ws.Range("I20").Value = 24.45785476576298
$ws.Range("K20").Value = 9.205636191829077
$ws.Range("L20").Value = 11.07301687562455
$ws.Range("M20").Value = 14.53351774484261
$ws.Range("N20").Value = 19.71482028658443
$ws.Range("O20").Value = 24.03876512177284

# Row 21
$ws.Range("B21").Value = 12.61923005530549
$ws.Range("C21").Value = 9.435901416260844
$ws.Range("D21").Value = 5.889176726429053
$ws.Range("F21").Value = 26.92259041010714
$ws.Range("G21").Value = 3.64604203177388
$ws.Range("I21").Value = 24.38061698677733
$ws.Range("K21").Value = 9.435866762547548
$ws.Range("L21").Value = 11.07635438801382
$ws.Range("M21").Value = 14.59663714577259
$ws.Range("N21").Value = 19.64750152054732
$ws.Range("O21").Value = 23.98549462711371

# Row 22
$ws.Range("B22").Value = 12.81796130084097
$ws.Range("C22").Value = 9.455711108327321
$ws.Range("D22").Value = 5.93403055981931
$ws.Range("F22").Value = 26.93233982301554
$ws.Range("G22").Value = 3.644876905920664
$ws.Range("I22").Value = 24.33371159494633
$ws.Range("K22").Value = 9.584273151955177
$ws.Range("L22").Value = 11.07998546928556
$ws.Range("M22").Value = 14.63955845516631
$ws.Range("N22").Value = 19.60493076164308
$ws.Range("O22").Value = 23.95464433147755

# Row 23
$ws.Range("B23").Value = 12.71214960351976
$ws.Range("C23").Value = 9.44515252154406
$ws.Range("D23").Value = 5.910179076805633
$ws.Range("F23").Value = 26.92664160683381
$ws.Range("G23").Value = 3.64549452952549
$ws.Range("I23").Value = 24.35841869462611
$ws.Range("K23").Value = 9.505270085257166
$ws.Range("L23").Value = 11.07791527126168
$ws.Range("M23").Value = 14.61650057932871
$ws.Range("N23").Value = 19.62752145678386
$ws.Range("O23").Value = 23.97074665141936

# Row 24
$ws.Range("B24").Value = 12.3061755354113
$ws.Range("C24").Value = 9.404858658112701
$ws.Range("D24").Value = 5.817946326455072
$ws.Range("F24").Value = 26.91594447987328
$ws.Range("G24").Value = 3.647926730764545
$ws.Range("I24").Value = 24.45917747402341
$ws.Range("K24").Value = 9.201841921018481
$ws.Range("L24").Value = 11.07298689567812
$ws.Range("M24").Value = 14.5325148755316
$ws.Range("N24").Value = 19.7159447292557
$ws.Range("O24").Value = 24.03970267649499

# Row 25
$ws.Range("B25").Value = 11.85957413722668
$ws.Range("C25").Value = 9.360791130796638
$ws.Range("D25").Value = 5.714603147397083
$ws.Range("F25").Value = 26.9277294327972
$ws.Range("G25").Value = 3.65075038684534
$ws.Range("I25").Value = 24.58318864327425
$ws.Range("K25").Value = 8.867342055486402
$ws.Range("L25").Value = 11.07393584701468
$ws.Range("M25").Value = 14.44935696434708
$ws.Range("N25").Value = 19.81748196813836
$ws.Range("O25").Value = 24.13107611901283
